# Add a "Spain" test-data sheet, cloned from the existing "Italy" sheet,
# and populate it with the Spain-market values (mirrors the Zettler Spain
# market entry pattern used by the other country sheets).

$wb = $excel.ActiveWorkbook

# The "Italy" sheet is the template for every other "<Country> Market" tab.
$italy = $wb.Worksheets.Item("Italy")

# Clear Italy's own selection down to a deliberate "select-all" state before
# we leave it (matches the author having selected the whole used range on
# the source sheet right before creating the new one).
$italy.Activate()
[void]$italy.Range("A1:D36").Select()

# Copy "Italy" and drop the clone immediately after it -> becomes sheet #8.
$italy.Copy($null, $italy)

$spain = $wb.ActiveSheet
$spain.Name = "Spain"

# Market name + identifier cells (same layout as every other country sheet).
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2019/T2046/T2066"

# Column B no longer needs to be as wide as "Italy Market" required; let
# Excel re-fit it (and the trailing notes column D) to the new content.
$spain.Columns.Item(2).AutoFit()
$spain.Columns.Item(4).AutoFit()

# Leave the new sheet active, with the cursor parked on E15 as the last
# touched cell.
$spain.Activate()
[void]$spain.Range("E15").Select()
